$d = $word.ActiveDocument
$CR = [char]13

# --- Step 1: merge the closing parenthesis into the first run's text ---
# "Associated salary management (HR accountant" -> "...accountant)"
$null = $d.Content.Find.Execute(
    "Associated salary management (HR accountant", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Associated salary management (HR accountant)", 2)

# Locate the paragraph that now reads "...accountant))" -- the original
# trailing ")" run (after the bookmark) is still present, duplicating the
# parenthesis just merged above. (Note: Paragraph.Index is unreliable here,
# so track the 1-based position in $d.Paragraphs ourselves.)
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.TrimEnd($CR)
    if ($text -eq "Associated salary management (HR accountant))") {
        $targetIdx = $i
        break
    }
}

$target = $d.Paragraphs.Item($targetIdx)
$paraEnd = $target.Range.End
# Character just before the paragraph mark is the stray ")" left over from
# the now-orphaned run; remove it, restoring a single trailing ")".
$extra = $d.Range($paraEnd - 2, $paraEnd - 1)
$extra.Delete()

# --- Step 2: insert a new list paragraph after it ---
$target = $d.Paragraphs.Item($targetIdx)
$target.Range.InsertParagraphAfter()

$newIdx = $targetIdx + 1
$newPara = $d.Paragraphs.Item($newIdx)

# Use a placeholder trailing character while placing the (collapsed)
# "_GoBack" bookmark: placing a zero-length bookmark exactly at a
# paragraph's end position is mishandled by the host, so we bookmark just
# before a sentinel character and then delete the sentinel, leaving the
# bookmark correctly collapsed at the real paragraph end.
$newPara.Range.Text = "Request for equipment/material#"
$newPara = $d.Paragraphs.Item($newIdx)

$sentinelStart = $newPara.Range.End - 2
$bmPoint = $d.Range($sentinelStart, $sentinelStart)
$d.Bookmarks.Add("_GoBack", $bmPoint)

$sentinel = $d.Range($sentinelStart, $sentinelStart + 1)
$sentinel.Delete()
